$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    if ($val -match '^-?[0-9]+(\.[0-9]+)?$') {
        $ws.Range($cell).Value = "'" + $val
    } else {
        $ws.Range($cell).Value = $val
    }
}

Set-TextValue 'D2' '26.101.16'
$ws.Range('E2').Value = '  -0.34%  '

Set-TextValue 'D3' '1.652.85'
$ws.Range('E3').Value = '  -0.74%  '

Set-TextValue 'D4' '1.002'
$ws.Range('E4').Value = '  -0.02%  '

Set-TextValue 'D5' '217.47'
$ws.Range('E5').Value = '  -0.17%  '

Set-TextValue 'D6' '0.5283'
$ws.Range('E6').Value = '  +0.44%  '

$ws.Range('E7').Value = '  +0.10%  '

Set-TextValue 'D8' '0.2606'
$ws.Range('E8').Value = '  -1.48%  '

Set-TextValue 'D9' '0.06316'
$ws.Range('E9').Value = '  +0.49%  '

Set-TextValue 'D10' '20.34'
$ws.Range('E10').Value = '  -2.12%  '

Set-TextValue 'D11' '0.07796'
$ws.Range('E11').Value = '  +0.53%  '

Set-TextValue 'D12' '4.520'
$ws.Range('E12').Value = '  +1.23%  '

Set-TextValue 'D13' '1.649.76'
$ws.Range('E13').Value = '  -0.56%  '

Set-TextValue 'D14' '1.879.61'
$ws.Range('E14').Value = '  -0.56%  '

Set-TextValue 'D15' '0.5478'
$ws.Range('E15').Value = '  -0.02%  '

Set-TextValue 'D16' '0.0₅8201'
$ws.Range('E16').Value = '  +0.94%  '

Set-TextValue 'D17' '65.36'
$ws.Range('E17').Value = '  +0.63%  '

Set-TextValue 'D18' '26.094.99'
$ws.Range('E18').Value = '  -0.29%  '

$ws.Range('E19').Value = '  -0.03%  '

Set-TextValue 'D20' '4.585'
$ws.Range('E20').Value = '  -0.25%  '

Set-TextValue 'D21' '190.69'
$ws.Range('E21').Value = '  -0.75%  '

$ws.Range('E22').Value = '  +0.23%  '

Set-TextValue 'D23' '6.012'
$ws.Range('E23').Value = '  -0.04%  '

$ws.Range('E24').Value = '  -0.04%  '

Set-TextValue 'D25' '144.73'
$ws.Range('E25').Value = '  +5.22%  '

Set-TextValue 'D26' '0.1228'
$ws.Range('E26').Value = '  -0.91%  '

Set-TextValue 'D27' '7.217'
$ws.Range('E27').Value = '  -0.56%  '

Set-TextValue 'D28' '15.98'
$ws.Range('E28').Value = '  -1.48%  '

Set-TextValue 'D29' '1.449'
$ws.Range('E29').Value = '  +3.36%  '

Set-TextValue 'D30' '0.05782'
$ws.Range('E30').Value = '  -3.34%  '

Set-TextValue 'D31' '1.272'
$ws.Range('E31').Value = '  -0.53%  '

Set-TextValue 'D32' '3.547'
$ws.Range('E32').Value = '  +0.44%  '

Set-TextValue 'D33' '3.262'
$ws.Range('E33').Value = '  -0.22%  '

Set-TextValue 'D34' '1.598'
$ws.Range('E34').Value = '  +1.08%  '

Set-TextValue 'D35' '2.795'
$ws.Range('E35').Value = '  +0.97%  '

Set-TextValue 'D36' '2.413'
$ws.Range('E36').Value = '  -0.14%  '

Set-TextValue 'D37' '0.9461'
$ws.Range('E37').Value = '  -1.61%  '

Set-TextValue 'D38' '0.5741'

Set-TextValue 'D39' '0.01612'
$ws.Range('E39').Value = '  +1.00%  '

Set-TextValue 'D40' '0.8520'
$ws.Range('E40').Value = '  +0.05%  '

$ws.Range('B41').Value = 'PaxDollar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue 'D41' '1.002'
$ws.Range('E41').Value = '  +0.06%  '

$ws.Range('B42').Value = 'Quant'
$ws.Range('C42').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue 'D42' '104.09'
$ws.Range('E42').Value = '  +3.11%  '

Set-TextValue 'D43' '5.711'
$ws.Range('E43').Value = '  -3.58%  '

Set-TextValue 'D44' '1.031.44'
$ws.Range('E44').Value = '  +2.28%  '

Set-TextValue 'D45' '1.793.95'
$ws.Range('E45').Value = '  -0.55%  '

$ws.Range('E46').Value = '  +0.08%  '

Set-TextValue 'D47' '1.000'
$ws.Range('E47').Value = '  -0.19%  '

Set-TextValue 'D48' '0.4330'
$ws.Range('E48').Value = '  -0.17%  '

Set-TextValue 'D49' '7.861'
$ws.Range('E49').Value = '  -1.95%  '

Set-TextValue 'D50' '0.05143'
$ws.Range('E50').Value = '  -0.09%  '

Set-TextValue 'D51' '1.447'
$ws.Range('E51').Value = '  -1.44%  '
